$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (H1, the last header cell) onto the
# two new header cells before writing their text, so I1/J1 pick up the same
# cell style (bold, bordered, centered) instead of creating a brand new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
  @(3,4), @(6,7), @(8,8), @(8,8), @(8,8), @(10,10), @(8,8), @(8,8), @(7,8), @(8,8),
  @(7,7), @(8,8), @(7,8), @(7,7), @(7,7), @(5,5), @(7,8), @(7,7), @(7,8), @(7,8),
  @(7,7), @(6,7), @(7,7), @(8,8), @(7,7), @(7,8), @(7,7), @(7,7), @(7,7), @(8,8),
  @(7,7), @(6,7), @(6,6), @(3,4), @(10,10), @(8,8), @(6,6), @(6,7), @(5,5), @(1,4),
  @(1,4), @(2,5), @(6,6), @(3,5), @(7,7), @(4,5), @(1,3), @(6,8), @(7,8), @(3,3),
  @(3,3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 9).Value = $data[$i][0]
  $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
